# "Updated on Negative Testing"
# Add a new negative test case (TC002 / Register) to the MasterTestCases
# sheet, styled to match the existing TC001 row, and mark it "Blocked"
# using a new bold/blue status style (mirrors the existing bold/green
# "Pass" style used for TC001's Status cell).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("MasterTestCases")
$ws2 = $wb.Worksheets.Item("ApplicationLogin")

# --- New row of test data -------------------------------------------------
$ws1.Range("A3").Value = "TC002"
$ws1.Range("B3").Value = "Register"
$ws1.Range("C3").Value = "N"
$ws1.Range("D3").Value = "Blocked"

# --- Borders: reuse the plain thin-box cell format already used on the
# "ApplicationLogin" sheet's body rows so the style is shared/deduped
# rather than creating a redundant duplicate style.
$ws2.Range("A3").Copy() | Out-Null
$ws1.Range("A2:C2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws1.Range("A3:C3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Status font: bold blue for the new "Blocked" result, same treatment
# the sheet already uses for bold green "Pass".
$d3 = $ws1.Range("D3")
$d3.Font.Bold = $true
$d3.Font.Color = 16711680   # RGB(0,0,255) -> blue, OLE BGR-packed value

# --- Restore cursor / active sheet positions -------------------------------
$ws1.Activate()
$ws1.Range("C4").Select() | Out-Null

$ws2.Activate()
$ws2.Range("B3").Select() | Out-Null
